# ---------------------------------------------------------------------------
# worknotes.xlsx edit:
#   - rename Sheet1 -> 注册
#   - add a new sheet 用户设置 after 登录
#   - rewrite the content of 登录 (client/server login sequence notes)
#   - populate 用户设置 (client/server user-info + account settings notes)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- workbook-level: rename Sheet1, add 用户设置 after 登录 ----------------
$wsReg = $wb.Worksheets.Item(1)
$wsReg.Name = "注册"

$wsLogin = $wb.Worksheets.Item("登录")
# Copy 登录 (rather than Worksheets.Add) so the new sheet inherits the same
# xml namespace set / sheetFormatPr as the rest of the workbook; its content
# gets fully overwritten below.
$wsLogin.Copy($null, $wsLogin)
$wsUser = $wb.Worksheets.Item(3)
$wsUser.Name = "用户设置"

# ---------------------------------------------------------------------------
# 登录 sheet: clear old content, write new client/server login flow notes
# ---------------------------------------------------------------------------
$wsLogin.Cells.Clear()

$wsLogin.Range("C1").Value = "client"
$wsLogin.Range("F1").Value = "server"

$wsLogin.Range("A4").Value = 43573
$wsLogin.Range("A4").NumberFormat = "mm-dd-yy"
$wsLogin.Range("C4").Value = "登录"
$wsLogin.Range("F4").Value = "controller - Users"

$wsLogin.Range("D5").Value = "api -store - Login"
$wsLogin.Range("G5").Value = "method - login"

$wsLogin.Range("H6").Value = "db api - login"

$wsLogin.Range("F7").Value = "reposonse client"

$wsLogin.Range("C8").Value = "收到response"
$wsLogin.Range("G8").Value = "{active, detailed_info_done}"

$wsLogin.Range("C9").Value = "成功"

$wsLogin.Range("D10").Value = "保存Token"

$wsLogin.Range("D11").Value = "store {active, detailed_info}"

$wsLogin.Range("D12").Value = "detailed_info_done = 0 跳转至“用户设置”页面"

$wsLogin.Columns.Item(1).ColumnWidth = 9.22
$wsLogin.Columns.Item(3).ColumnWidth = 13.65
$wsLogin.Columns.Item(4).ColumnWidth = 29.36

$wsLogin.Range("D16").Select()

# ---------------------------------------------------------------------------
# 用户设置 sheet: user-info tab + account-management tab flow notes
# ---------------------------------------------------------------------------
$wsUser.Cells.Clear()

$wsUser.Range("C1").Value = "client"
$wsUser.Range("F1").Value = "server"

$wsUser.Range("A4").Value = 43573
# Re-use the date cell style already created on 登录!A4 (paste formats only)
# instead of calling NumberFormat again, which would otherwise create a
# duplicate (but equivalent) cellXfs entry.
$wsLogin.Range("A4").Copy()
$wsUser.Range("A4").PasteSpecial(-4122)
$wsUser.Range("C4").Value = "向server request用户信息，页面等待加载"
$wsUser.Range("F4").Value = "controller - Users"

$wsUser.Range("D5").Value = "api - requestUserInfo"
$wsUser.Range("G5").Value = "method - request_user_info"

$wsUser.Range("H6").Value = "db api - get_user_info"

$wsUser.Range("F7").Value = "reposonse client"

$wsUser.Range("C8").Value = "收到response，取消页面加载"

$wsUser.Range("C11").Value = "用户信息Tab"
$wsUser.Range("F11").Value = "controller - Users"

$wsUser.Range("D12").Value = "api - updateUserInfo"
$wsUser.Range("G12").Value = "method - update_user_info"

$wsUser.Range("H13").Value = "detailed_info_done 置1"

$wsUser.Range("C15").Value = "账号管理Tab"

$wsUser.Range("D16").Value = "api - updateEmail"
$wsUser.Range("D17").Value = "api - updatePassword"
$wsUser.Range("D18").Value = "api - updatePhone"

$wsUser.Columns.Item(1).ColumnWidth = 9.22
$wsUser.Columns.Item(4).ColumnWidth = 22.93
$wsUser.Columns.Item(6).ColumnWidth = 18.65
$wsUser.Columns.Item(7).ColumnWidth = 20.36
$wsUser.Columns.Item(8).ColumnWidth = 19.65

$wsUser.PageSetup.PaperSize = 9
$wsUser.PageSetup.Orientation = 1

$wsUser.Range("E18").Select()
$wsUser.Tab.Selected = $false

# --- restore 登录 as the active/selected sheet (unchanged from before) -----
$wsLogin.Activate()
